$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("D9 ->")

# Update "CURRENT AS OF:" label from D11.1 to D11.2
$ws.Range("V1").Value = "D11.2"

# Row 8: weather - J8 value, V8 total (font turns red as value decreased)
$ws.Range("J8").Value = 2
$ws.Range("V8").Font.ColorIndex = 3

# Row 9: time of day - J9 value, V9 total (font turns red as value decreased)
$ws.Range("J9").Value = 4
$ws.Range("V9").Font.ColorIndex = 3

# Row 10: total highlighted red too (value unchanged, only formatting)
$ws.Range("V10").Font.ColorIndex = 3

# Ordnance updates
$ws.Range("J12").Value = 3
$ws.Range("J15").Value = 9
$ws.Range("J18").Value = 9

# Update active selection to K12
$ws.Activate()
$ws.Range("K12").Select()
